# "adding area to discharge files"
#
# Adds two new columns to the discharge-computation sheet:
#   G: "Area"   -> per-segment cross-sectional area (D{r}-D{r-1})*B{r}/100
#   H: "Atotal" -> total area = SUM(G2:G11)
#
# Row 2 is the special case: there's no "previous" depth row, so the area
# formula treats the previous depth as 0: (D2-0)*B2/100.
# Rows 4:11 share one formula (relative refs), matching how Excel collapses
# a fill-down into a single shared formula group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (row 1)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area column (G): per-row incremental area
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

# Total area (H2)
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Matches the saved selection in the edited workbook
$ws.Range("E9").Select()
